$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf18"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.4765323333333333
$ws.Range("H2").Value = 1.429597
$ws.Range("I2").Value = 0.04733670697480491
$ws.Range("J2").Value = 0.04733670697480491
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.656156333333334
$ws.Range("N2").Value = 4.968469000000001
$ws.Range("O2").Value = 0.6151212440816572
$ws.Range("P2").Value = 0.6151212440816572
$ws.Range("Q2").Value = 0.7892120418881112
$ws.Range("R2").Value = 7.102908376993001
$ws.Range("S2").Value = 0.02911781408507086
$ws.Range("T2").Value = 0.02911781408507086

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf18"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.4765323333333333
$ws.Range("H3").Value = 1.429597
$ws.Range("I3").Value = 0.04733670697480491
$ws.Range("J3").Value = 0.04733670697480491
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.5648773333333333
$ws.Range("N3").Value = 1.694632
$ws.Range("O3").Value = 0.2098038941373262
$ws.Range("P3").Value = 0.2098038941373262
$ws.Range("Q3").Value = 0.2691823137004444
$ws.Range("R3").Value = 2.422640823304
$ws.Range("S3").Value = 0.0099314254589516
$ws.Range("T3").Value = 0.009931425458951598

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf18"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.4765323333333333
$ws.Range("H4").Value = 1.429597
$ws.Range("I4").Value = 0.04733670697480491
$ws.Range("J4").Value = 0.04733670697480491
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4713726666666667
$ws.Range("N4").Value = 1.414118
$ws.Range("O4").Value = 0.1750748617810164
$ws.Range("P4").Value = 0.1750748617810165
$ws.Range("Q4").Value = 0.2246243167162222
$ws.Range("R4").Value = 2.021618850446
$ws.Range("S4").Value = 0.008287467430782446
$ws.Range("T4").Value = 0.008287467430782448

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf18"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.770835333333332
$ws.Range("H5").Value = 26.312506
$ws.Range("I5").Value = 0.8712576945074703
$ws.Range("J5").Value = 0.8712576945074703
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.656156333333334
$ws.Range("N5").Value = 4.968469000000001
$ws.Range("O5").Value = 0.6151212440816572
$ws.Range("P5").Value = 0.6151212440816572
$ws.Range("Q5").Value = 14.52587448592378
$ws.Range("R5").Value = 130.732870373314
$ws.Range("S5").Value = 0.5359291169611515
$ws.Range("T5").Value = 0.5359291169611515

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf18"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 8.770835333333332
$ws.Range("H6").Value = 26.312506
$ws.Range("I6").Value = 0.8712576945074703
$ws.Range("J6").Value = 0.8712576945074703
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.5648773333333333
$ws.Range("N6").Value = 1.694632
$ws.Range("O6").Value = 0.2098038941373262
$ws.Range("P6").Value = 0.2098038941373262
$ws.Range("Q6").Value = 4.954446074199111
$ws.Range("R6").Value = 44.590014667792
$ws.Range("S6").Value = 0.1827932571047762
$ws.Range("T6").Value = 0.1827932571047761

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf18"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 8.770835333333332
$ws.Range("H7").Value = 26.312506
$ws.Range("I7").Value = 0.8712576945074703
$ws.Range("J7").Value = 0.8712576945074703
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4713726666666667
$ws.Range("N7").Value = 1.414118
$ws.Range("O7").Value = 0.1750748617810164
$ws.Range("P7").Value = 0.1750748617810165
$ws.Range("Q7").Value = 4.134332039967555
$ws.Range("R7").Value = 37.208988359708
$ws.Range("S7").Value = 0.1525353204415424
$ws.Range("T7").Value = 0.1525353204415424

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fgf18"
$ws.Range("C8").Value = "Fgfr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8194993333333332
$ws.Range("H8").Value = 2.458498
$ws.Range("I8").Value = 0.08140559851772487
$ws.Range("J8").Value = 0.08140559851772486
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.656156333333334
$ws.Range("N8").Value = 4.968469000000001
$ws.Range("O8").Value = 0.6151212440816572
$ws.Range("P8").Value = 0.6151212440816572
$ws.Range("Q8").Value = 1.357219011062444
$ws.Range("R8").Value = 12.214971099562
$ws.Range("S8").Value = 0.05007431303543484
$ws.Range("T8").Value = 0.05007431303543483

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fgf18"
$ws.Range("C9").Value = "Fgfr3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8194993333333332
$ws.Range("H9").Value = 2.458498
$ws.Range("I9").Value = 0.08140559851772487
$ws.Range("J9").Value = 0.08140559851772486
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.5648773333333333
$ws.Range("N9").Value = 1.694632
$ws.Range("O9").Value = 0.2098038941373262
$ws.Range("P9").Value = 0.2098038941373262
$ws.Range("Q9").Value = 0.4629165980817778
$ws.Range("R9").Value = 4.166249382735999
$ws.Range("S9").Value = 0.01707921157359843
$ws.Range("T9").Value = 0.01707921157359842

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fgf18"
$ws.Range("C10").Value = "Fgfr3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8194993333333332
$ws.Range("H10").Value = 2.458498
$ws.Range("I10").Value = 0.08140559851772487
$ws.Range("J10").Value = 0.08140559851772486
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4713726666666667
$ws.Range("N10").Value = 1.414118
$ws.Range("O10").Value = 0.1750748617810164
$ws.Range("P10").Value = 0.1750748617810165
$ws.Range("Q10").Value = 0.3862895860848889
$ws.Range("R10").Value = 3.476606274764
$ws.Range("S10").Value = 0.0142520739086916
$ws.Range("T10").Value = 0.0142520739086916

